$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking strings (e.g. "14.80", "0.9993") are not
# auto-coerced to numbers by Excel, matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.881.07'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.802.76'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '309.65'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').Value = '0.4398'
$ws.Range('E7').Value = '  +4.21%  '
$ws.Range('D8').Value = '0.3685'
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('D9').Value = '0.07425'
$ws.Range('E9').Value = '  +2.86%  '
$ws.Range('D10').Value = '0.8591'
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('D11').Value = '20.77'
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').Value = '1.799.22'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').Value = '6.620'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').Value = '93.08'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').Value = '0.07068'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '5.271'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '0.000008677'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '0.9995'
$ws.Range('D20').Value = '14.80'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').Value = '26.893.57'
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').Value = '5.158'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('D23').Value = '10.81'
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = '1.971'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = '151.70'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '2.211'
$ws.Range('E26').Value = '  -1.82%  '
$ws.Range('D27').Value = '18.42'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').Value = '5.192'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('D29').Value = '117.60'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '0.08795'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').Value = '0.7419'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').Value = '1.162'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('D33').Value = '4.466'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').Value = '0.9989'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '0.01964'
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('D38').Value = '0.05199'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('D39').Value = '0.5253'
$ws.Range('E39').Value = '  +4.19%  '
$ws.Range('D40').Value = '7.052'
$ws.Range('E40').Value = '  -3.64%  '
$ws.Range('D41').Value = '2.812'
$ws.Range('E41').Value = '  -1.98%  '
$ws.Range('D42').Value = '0.1680'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = '8.469'
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('D44').Value = '0.4981'
$ws.Range('E44').Value = '  +5.88%  '
$ws.Range('D45').Value = '2.047'
$ws.Range('E45').Value = '  +7.06%  '
$ws.Range('D46').Value = '10.38'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').Value = '104.35'
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('D48').Value = '0.9988'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').Value = '1.667'
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('D50').Value = '0.06313'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').Value = '0.9202'
$ws.Range('E51').Value = '  +2.41%  '

# Remove the temporary text-number-format so column D keeps its original (default) style.
$ws.Range("D2:D51").ClearFormats()
